$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.920.52"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "1.856.65"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'316.96"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.4353"
$ws.Range("E7").Value = "  -5.27%  "
$ws.Range("D8").Value = "'0.3676"
$ws.Range("E8").Value = "  -3.80%  "
$ws.Range("D9").Value = "'0.07487"
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("D10").Value = "'0.9378"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").Value = "'21.35"
$ws.Range("E11").Value = "  -3.34%  "
$ws.Range("D12").Value = "1.830.74"
$ws.Range("E12").Value = "  -4.56%  "
$ws.Range("D13").Value = "'6.704"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("D14").Value = "'5.426"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").Value = "'0.06865"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'81.35"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "'0.000009025"
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'15.87"
$ws.Range("E20").Value = "  -4.95%  "
$ws.Range("D21").Value = "27.915.23"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").Value = "'5.097"
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").Value = "'11.02"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "2.073.91"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "'2.008"
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("D26").Value = "'154.22"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").Value = "'18.33"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").Value = "'5.405"
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").Value = "'113.25"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("D30").Value = "'1.730"
$ws.Range("E30").Value = "  -6.76%  "
$ws.Range("D31").Value = "'0.08955"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").Value = "'0.8026"
$ws.Range("E32").Value = "  -7.58%  "
$ws.Range("D33").Value = "'4.818"
$ws.Range("E33").Value = "  -5.21%  "
$ws.Range("D34").Value = "'3.003"
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("D35").Value = "'1.167"
$ws.Range("E35").Value = "  -6.75%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'1.116"
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("D39").Value = "'0.01966"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").Value = "'2.924"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D42").Value = "'7.016"
$ws.Range("E42").Value = "  -5.24%  "
$ws.Range("D43").Value = "'0.1680"
$ws.Range("E43").Value = "  -4.19%  "
$ws.Range("D44").Value = "'8.737"
$ws.Range("E44").Value = "  -6.09%  "
$ws.Range("D45").Value = "'0.06713"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").Value = "'0.4873"
$ws.Range("E46").Value = "  -5.83%  "
$ws.Range("D47").Value = "'10.63"
$ws.Range("E47").Value = "  -5.28%  "
$ws.Range("D48").Value = "'106.73"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("D49").Value = "'1.927"
$ws.Range("E49").Value = "  -7.85%  "
$ws.Range("E50").Value = "  -5.87%  "
$ws.Range("E51").Value = "  -0.17%  "
